$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final "Out of PO" table data (header + 17 players), replacing the prior
# 18-player list (Jabari Smith Jr. removed, other rows refreshed/reordered).
$data = @(
    @("Oyuncu Adı", "Pozisyon", "Takım"),
    @("Payton Pritchard", "PG", "Boston Celtics"),
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Deandre Ayton", "C", "Portland Trail Blazers"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Robert Williams III", "C", "Portland Trail Blazers"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Tari Eason", "SF,PF", "Houston Rockets")
)

$oldLastRow = 19
$newLastRow = $data.Count

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# The source list shrank by one row (18 -> 17 players); clear the now-unused
# trailing row so it doesn't linger with stale values.
if ($newLastRow -lt $oldLastRow) {
    $clearRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($oldLastRow, 3))
    $clearRange.Clear()
}
